$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H19").Value = 2302.625
$ws_ALC.Range("I19").Value = 2268
$ws_ALC.Range("K19").Value = 2268
$ws_ALC.Range("M19").Value = -2093
$ws_ALC.Range("H42").Value = 3369.6667
$ws_ALC.Range("J42").Value = 6687.6665
$ws_ALC.Range("L42").Value = 20062.9995
$ws_ALC.Range("N42").Value = -20522.9995
$ws_ALC.Range("H69").Value = 0
$ws_ALC.Range("J69").Value = 0
$ws_ALC.Range("L69").Value = 0
$ws_ALC.Range("N69").ClearContents()
$ws_ALC.Range("H72").Value = 0
$ws_ALC.Range("J72").Value = 0
$ws_ALC.Range("L72").Value = 0
$ws_ALC.Range("N72").ClearContents()
$ws_ALC.Range("H101").Value = 577.6
$ws_ALC.Range("I101").Value = 577.6
$ws_ALC.Range("K101").Value = 1732.8
$ws_ALC.Range("M101").Value = -110.8000000000002
$ws_ALC.Range("H107").Value = 137.57143
$ws_ALC.Range("I107").Value = 138.8
$ws_ALC.Range("J107").Value = 134.5
$ws_ALC.Range("K107").Value = 138.8
$ws_ALC.Range("L107").Value = 134.5
$ws_ALC.Range("M107").Value = 1781.2
$ws_ALC.Range("N107").Value = -3974.5
$ws_ALC.Range("H113").Value = 16599.666
$ws_ALC.Range("I113").Value = 13949.5
$ws_ALC.Range("K113").Value = 13949.5
$ws_ALC.Range("M113").Value = -10695.5
$ws_ALC.Range("H137").Value = 1617.8182
$ws_ALC.Range("I137").Value = 709.6
$ws_ALC.Range("J137").Value = 2374.6667
$ws_ALC.Range("K137").Value = 2128.8
$ws_ALC.Range("L137").Value = 7124.000100000001
$ws_ALC.Range("M137").Value = 421.1999999999998
$ws_ALC.Range("N137").Value = -12224.0001
$ws_ALC.Range("H138").Value = 3535.5483
$ws_ALC.Range("I138").Value = 2655.3333
$ws_ALC.Range("J138").Value = 3629.8572
$ws_ALC.Range("K138").Value = 7965.999899999999
$ws_ALC.Range("L138").Value = 10889.5716
$ws_ALC.Range("M138").Value = -2825.999899999999
$ws_ALC.Range("N138").Value = -21169.5716
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H74").Value = 3059.2
$ws_ARM.Range("I74").Value = 3059.2
$ws_ARM.Range("K74").Value = 3059.2
$ws_ARM.Range("M74").Value = -2185.2
$ws_ARM.Range("H77").Value = 3059.2
$ws_ARM.Range("I77").Value = 3059.2
$ws_ARM.Range("K77").Value = 15296
$ws_ARM.Range("M77").Value = -10928
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H86").Value = 2463.8333
$ws_BSM.Range("I86").Value = 2033
$ws_BSM.Range("K86").Value = 2033
$ws_BSM.Range("M86").Value = -910
$ws_BSM.Range("H89").Value = 2463.8333
$ws_BSM.Range("I89").Value = 2033
$ws_BSM.Range("K89").Value = 10165
$ws_BSM.Range("M89").Value = -4549
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H5").Value = 193.33333
$ws_CRP.Range("I5").Value = 217.5
$ws_CRP.Range("K5").Value = 217.5
$ws_CRP.Range("M5").Value = -105.5
$ws_CRP.Range("H107").Value = 1281.125
$ws_CRP.Range("I107").Value = 916.3333
$ws_CRP.Range("K107").Value = 916.3333
$ws_CRP.Range("M107").Value = 1003.6667
$ws_CRP.Range("H132").Value = 1585.6
$ws_CRP.Range("I132").Value = 1372.8889
$ws_CRP.Range("K132").Value = 4118.6667
$ws_CRP.Range("M132").Value = -1588.6667
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H17").Value = 1007.1667
$ws_CUL.Range("I17").Value = 23.5
$ws_CUL.Range("J17").Value = 1499
$ws_CUL.Range("K17").Value = 70.5
$ws_CUL.Range("L17").Value = 4497
$ws_CUL.Range("M17").Value = 98.5
$ws_CUL.Range("N17").Value = -4835
$ws_CUL.Range("H19").Value = 0
$ws_CUL.Range("J19").Value = 0
$ws_CUL.Range("L19").Value = 0
$ws_CUL.Range("N19").ClearContents()
$ws_CUL.Range("H22").Value = 870
$ws_CUL.Range("I22").Value = 40
$ws_CUL.Range("J22").Value = 1700
$ws_CUL.Range("K22").Value = 120
$ws_CUL.Range("L22").Value = 5100
$ws_CUL.Range("N22").Value = -5438
$ws_CUL.Range("M22").Value = 49
$ws_CUL.Range("H26").Value = 105723.1
$ws_CUL.Range("J26").Value = 10284.2
$ws_CUL.Range("L26").Value = 30852.6
$ws_CUL.Range("N26").Value = -31428.6
$ws_CUL.Range("H27").Value = 870
$ws_CUL.Range("I27").Value = 40
$ws_CUL.Range("J27").Value = 1700
$ws_CUL.Range("K27").Value = 120
$ws_CUL.Range("L27").Value = 5100
$ws_CUL.Range("N27").Value = -5304
$ws_CUL.Range("M27").Value = -18
$ws_CUL.Range("H37").Value = 99939
$ws_CUL.Range("J37").Value = 99939
$ws_CUL.Range("L37").Value = 299817
$ws_CUL.Range("N37").Value = -300041
$ws_CUL.Range("H128").Value = 340295.38
$ws_CUL.Range("I128").Value = 340295.38
$ws_CUL.Range("K128").Value = 1020886.14
$ws_CUL.Range("M128").Value = -1015906.14
$ws_CUL.Range("H130").Value = 2865.4
$ws_CUL.Range("J130").Value = 3698.5
$ws_CUL.Range("L130").Value = 11095.5
$ws_CUL.Range("N130").Value = -21135.5
$ws_CUL.Range("H131").Value = 925
$ws_CUL.Range("J131").Value = 0
$ws_CUL.Range("L131").Value = 0
$ws_CUL.Range("N131").ClearContents()
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H25").Value = 2400
$ws_GSM.Range("J25").Value = 800
$ws_GSM.Range("L25").Value = 800
$ws_GSM.Range("N25").Value = -1858
$ws_GSM.Range("H80").Value = 3411.3333
$ws_GSM.Range("I80").Value = 3245
$ws_GSM.Range("K80").Value = 3245
$ws_GSM.Range("M80").Value = -2247
$ws_GSM.Range("H83").Value = 3411.3333
$ws_GSM.Range("I83").Value = 3245
$ws_GSM.Range("K83").Value = 16225
$ws_GSM.Range("M83").Value = -11233
$ws_GSM.Range("H97").Value = 1000
$ws_GSM.Range("I97").Value = 1000
$ws_GSM.Range("K97").Value = 1000
$ws_GSM.Range("M97").Value = -504
$ws_GSM.Range("H132").Value = 0
$ws_GSM.Range("I132").Value = 0
$ws_GSM.Range("K132").Value = 0
$ws_GSM.Range("M132").ClearContents()
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 1418.3636
$ws_LTW.Range("I22").Value = 1387.5
$ws_LTW.Range("J22").Value = 1500.6666
$ws_LTW.Range("K22").Value = 1387.5
$ws_LTW.Range("L22").Value = 1500.6666
$ws_LTW.Range("M22").Value = -1092.5
$ws_LTW.Range("N22").Value = -2090.6666
$ws_LTW.Range("H27").Value = 1418.3636
$ws_LTW.Range("I27").Value = 1387.5
$ws_LTW.Range("J27").Value = 1500.6666
$ws_LTW.Range("K27").Value = 1387.5
$ws_LTW.Range("L27").Value = 1500.6666
$ws_LTW.Range("M27").Value = -1280.5
$ws_LTW.Range("N27").Value = -1714.6666
$ws_LTW.Range("I43").Value = 0
$ws_LTW.Range("K43").Value = 0
$ws_LTW.Range("M43").ClearContents()
$ws_LTW.Range("H46").Value = 3261.6553
$ws_LTW.Range("I46").Value = 2399.3572
$ws_LTW.Range("K46").Value = 2399.3572
$ws_LTW.Range("M46").Value = -2211.3572
$ws_LTW.Range("H55").Value = 313.33334
$ws_LTW.Range("I55").Value = 278.75
$ws_LTW.Range("J55").Value = 330.625
$ws_LTW.Range("K55").Value = 278.75
$ws_LTW.Range("L55").Value = 330.625
$ws_LTW.Range("M55").Value = -105.75
$ws_LTW.Range("N55").Value = -676.625
